$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was inserted into the daily-logic subset sheet right
# before the existing row 263, pushing all the subsequent rows (old 263-307)
# down by one (new 264-308). Insert a fresh row 263 and shift everything else.
$ws.Rows.Item(263).EntireRow.Insert()

# Populate the newly inserted row 263 with the new market record.
$ws.Cells.Item(263, 1).Value = 4
$ws.Cells.Item(263, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(263, 3).Value = "Los Lagos"
$ws.Cells.Item(263, 4).Value = 44798
$ws.Cells.Item(263, 5).Value = 10
$ws.Cells.Item(263, 6).Value = 100112017
$ws.Cells.Item(263, 7).Value = "Apio"
$ws.Cells.Item(263, 8).Value = "Americana (o)"
$ws.Cells.Item(263, 9).Value = "Segunda"
$ws.Cells.Item(263, 10).Value = 30
$ws.Cells.Item(263, 11).Value = 12000
$ws.Cells.Item(263, 12).Value = 12000
$ws.Cells.Item(263, 13).Value = 12000
$ws.Cells.Item(263, 14).Value = "`$/docena de matas"
$ws.Cells.Item(263, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(263, 16).Value = 2000
$ws.Cells.Item(263, 17).Value = 6
$ws.Cells.Item(263, 18).Value = "Hortaliza"
